$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H28").Value = 741.8889
$ws.Range("I28").Value = 335.6
$ws.Range("K28").Value = 335.6
$ws.Range("M28").Value = 149.4
$ws.Range("H31").Value = 3463.5
$ws.Range("I31").Value = 2701.1428
$ws.Range("K31").Value = 8103.428400000001
$ws.Range("M31").Value = -7873.428400000001
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H58").Value = 807.5
$ws.Range("I58").Value = 807.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2422.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2272.5
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 14400
$ws.Range("I61").Value = 100
$ws.Range("J61").Value = 28700
$ws.Range("K61").Value = 300
$ws.Range("L61").Value = 86100
$ws.Range("M61").Value = -128
$ws.Range("N61").Value = -86444
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H138").Value = 2499.6
$ws.Range("J138").Value = 2749.75
$ws.Range("L138").Value = 8249.25
$ws.Range("N138").Value = -18529.25
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2500
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814
$ws.Range("H66").Value = 2500
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068
$ws.Range("H102").Value = 2249.5
$ws.Range("I102").Value = 2249.5
$ws.Range("K102").Value = 2249.5
$ws.Range("M102").Value = -627.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2550.0833
$ws.Range("I20").Value = 2540.25
$ws.Range("J20").Value = 2569.75
$ws.Range("K20").Value = 2540.25
$ws.Range("L20").Value = 2569.75
$ws.Range("M20").Value = -2293.25
$ws.Range("N20").Value = -3063.75
$ws.Range("H22").Value = 428.2857
$ws.Range("I22").Value = 432.66666
$ws.Range("J22").Value = 402
$ws.Range("K22").Value = 432.66666
$ws.Range("L22").Value = 402
$ws.Range("M22").Value = -259.66666
$ws.Range("N22").Value = -748
$ws.Range("H86").Value = 2111.9
$ws.Range("I86").Value = 2111.9
$ws.Range("K86").Value = 2111.9
$ws.Range("M86").Value = -988.9000000000001
$ws.Range("H89").Value = 2111.9
$ws.Range("I89").Value = 2111.9
$ws.Range("K89").Value = 10559.5
$ws.Range("M89").Value = -4943.5
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H107").Value = 1873.1666
$ws.Range("I107").Value = 750
$ws.Range("J107").Value = 2097.8
$ws.Range("K107").Value = 750
$ws.Range("L107").Value = 2097.8
$ws.Range("M107").Value = 1170
$ws.Range("N107").Value = -5937.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 126.5
$ws.Range("J7").Value = 214
$ws.Range("L7").Value = 214
$ws.Range("N7").Value = -440
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H68").Value = 49999
$ws.Range("J68").Value = 49999
$ws.Range("L68").Value = 49999
$ws.Range("N68").Value = -51497
$ws.Range("H71").Value = 49999
$ws.Range("J71").Value = 49999
$ws.Range("L71").Value = 149997
$ws.Range("N71").Value = -157485
$ws.Range("H74").Value = 39975
$ws.Range("J74").Value = 39975
$ws.Range("L74").Value = 39975
$ws.Range("N74").Value = -41723
$ws.Range("H77").Value = 39975
$ws.Range("J77").Value = 39975
$ws.Range("L77").Value = 119925
$ws.Range("N77").Value = -128661
$ws.Range("H99").Value = 5466.6665
$ws.Range("I99").Value = 4500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -3002
$ws.Range("H105").Value = 3489.8
$ws.Range("I105").Value = 4062.25
$ws.Range("K105").Value = 4062.25
$ws.Range("M105").Value = -2315.25
$ws.Range("H125").Value = 55000
$ws.Range("J125").Value = 55000
$ws.Range("L125").Value = 55000
$ws.Range("N125").Value = -59920
$ws.Range("H126").Value = 5466.6665
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030
$ws.Range("H134").Value = 7000
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = -465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 142867.14
$ws.Range("I6").Value = 16.5
$ws.Range("J6").Value = 333334.66
$ws.Range("K6").Value = 49.5
$ws.Range("L6").Value = 1000003.98
$ws.Range("M6").Value = 63.5
$ws.Range("N6").Value = -1000229.98
$ws.Range("H44").Value = 491.57144
$ws.Range("J44").Value = 491.57144
$ws.Range("L44").Value = 1474.71432
$ws.Range("N44").Value = -2270.71432
$ws.Range("H68").Value = 613.5
$ws.Range("J68").Value = 467.5
$ws.Range("L68").Value = 1402.5
$ws.Range("N68").Value = -3024.5
$ws.Range("H71").Value = 613.5
$ws.Range("J71").Value = 467.5
$ws.Range("L71").Value = 4207.5
$ws.Range("N71").Value = -12319.5
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H140").Value = 2150.25
$ws.Range("I140").Value = 2150.25
$ws.Range("K140").Value = 6450.75
$ws.Range("M140").Value = -1270.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 14038769
$ws.Range("I11").Value = 22062874
$ws.Range("J11").Value = 1200200.8
$ws.Range("K11").Value = 22062874
$ws.Range("L11").Value = 1200200.8
$ws.Range("M11").Value = -22062735
$ws.Range("N11").Value = -1200478.8
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 64999.668
$ws.Range("J64").Value = 64999.668
$ws.Range("L64").Value = 64999.668
$ws.Range("N64").Value = -65449.668
$ws.Range("H67").Value = 64999.668
$ws.Range("J67").Value = 64999.668
$ws.Range("L67").Value = 64999.668
$ws.Range("N67").Value = -66559.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 35000
$ws.Range("I2").Value = 35000
$ws.Range("K2").Value = 35000
$ws.Range("M2").Value = -34888
$ws.Range("H28").Value = 175009.5
$ws.Range("J28").Value = 175009.5
$ws.Range("L28").Value = 175009.5
$ws.Range("N28").Value = -175705.5
$ws.Range("H63").Value = 10000
$ws.Range("I63").Value = 10000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 10000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -9376
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 10000
$ws.Range("I66").Value = 10000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 30000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -26880
$ws.Range("N66").ClearContents()
